$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 326-328. This pushes the existing rows
# 326-329 (the "Especial/Primera/Segunda/Tercera" set dated 44335) down
# to rows 329-332, unchanged, and grows the sheet dimension accordingly.
$ws.Range("A326:A328").EntireRow.Insert()

# Populate the newly inserted rows with the new week's data (date 44595),
# reusing the same layout/formatting as the rows that were pushed down.

# Row 326: Primera
$ws.Cells.Item(326, 1).Value = 2
$ws.Cells.Item(326, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(326, 3).Value = "Coquimbo"
$ws.Cells.Item(326, 4).Value = 44595
$ws.Cells.Item(326, 5).Value = 4
$ws.Cells.Item(326, 6).Value = 100112043
$ws.Cells.Item(326, 7).Value = "Pepino dulce"
$ws.Cells.Item(326, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 400
$ws.Cells.Item(326, 11).Value = 13500
$ws.Cells.Item(326, 12).Value = 14000
$ws.Cells.Item(326, 13).Value = 13750
$ws.Cells.Item(326, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(326, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(326, 16).Value = 764
$ws.Cells.Item(326, 17).Value = 18
$ws.Cells.Item(326, 18).Value = "Hortaliza"

# Row 327: Segunda
$ws.Cells.Item(327, 1).Value = 2
$ws.Cells.Item(327, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(327, 3).Value = "Coquimbo"
$ws.Cells.Item(327, 4).Value = 44595
$ws.Cells.Item(327, 5).Value = 4
$ws.Cells.Item(327, 6).Value = 100112043
$ws.Cells.Item(327, 7).Value = "Pepino dulce"
$ws.Cells.Item(327, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(327, 9).Value = "Segunda"
$ws.Cells.Item(327, 10).Value = 400
$ws.Cells.Item(327, 11).Value = 11500
$ws.Cells.Item(327, 12).Value = 12000
$ws.Cells.Item(327, 13).Value = 11750
$ws.Cells.Item(327, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(327, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(327, 16).Value = 653
$ws.Cells.Item(327, 17).Value = 18
$ws.Cells.Item(327, 18).Value = "Hortaliza"

# Row 328: Tercera
$ws.Cells.Item(328, 1).Value = 2
$ws.Cells.Item(328, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(328, 3).Value = "Coquimbo"
$ws.Cells.Item(328, 4).Value = 44595
$ws.Cells.Item(328, 5).Value = 4
$ws.Cells.Item(328, 6).Value = 100112043
$ws.Cells.Item(328, 7).Value = "Pepino dulce"
$ws.Cells.Item(328, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(328, 9).Value = "Tercera"
$ws.Cells.Item(328, 10).Value = 360
$ws.Cells.Item(328, 11).Value = 8500
$ws.Cells.Item(328, 12).Value = 9000
$ws.Cells.Item(328, 13).Value = 8750
$ws.Cells.Item(328, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(328, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(328, 16).Value = 486
$ws.Cells.Item(328, 17).Value = 18
$ws.Cells.Item(328, 18).Value = "Hortaliza"

# Column D (Fecha) keeps the same date style/number format as the rest of
# the column (style index carried over automatically by Insert(), but set
# explicitly too for safety).
$ws.Range("D326:D328").NumberFormat = $ws.Range("D325").NumberFormat
